$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'60.250.05"
$c.Style = "Normal"
$ws.Range("E2").Value = '  +1.75%  '
$c = $ws.Range("D3")
$c.Value = "'2.607.74"
$c.Style = "Normal"
$ws.Range("E3").Value = '  +0.94%  '
$ws.Range("E4").Value = '  +0.12%  '
$c = $ws.Range("D5")
$c.Value = "'561.85"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.15%  '
$c = $ws.Range("D6")
$c.Value = "'142.35"
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.73%  '
$ws.Range("E7").Value = '  -0.16%  '
$c = $ws.Range("D8")
$c.Value = "'0.601"
$c.Style = "Normal"
$ws.Range("E8").Value = '  +0.29%  '
$c = $ws.Range("D9")
$c.Value = "'2.640.71"
$c.Style = "Normal"
$ws.Range("E9").Value = '  +2.03%  '
$c = $ws.Range("D10")
$c.Value = "'6.69"
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.55%  '
$c = $ws.Range("D11")
$c.Value = "'0.105"
$c.Style = "Normal"
$ws.Range("E11").Value = '  +1.47%  '
$ws.Range("E12").Value = '  +3.58%  '
$c = $ws.Range("D13")
$c.Value = "'0.371"
$c.Style = "Normal"
$ws.Range("E13").Value = '  +8.17%  '
$c = $ws.Range("D14")
$c.Value = "'3.080.21"
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.38%  '
$c = $ws.Range("D15")
$c.Value = "'60.253.13"
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.70%  '
$c = $ws.Range("D16")
$c.Value = "'23.36"
$c.Style = "Normal"
$ws.Range("E16").Value = '  +3.78%  '
$ws.Range("E17").Value = '  +1.33%  '
$c = $ws.Range("D18")
$c.Value = "'2.621.83"
$c.Style = "Normal"
$ws.Range("E18").Value = '  +1.34%  '
$ws.Range("E19").Value = '  +3.06%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range("D20")
$c.Value = "'344.35"
$c.Style = "Normal"
$ws.Range("E20").Value = '  +2.57%  '
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range("D21")
$c.Value = "'10.81"
$c.Style = "Normal"
$ws.Range("E21").Value = '  +5.90%  '
$c = $ws.Range("D22")
$c.Value = "'6.89"
$c.Style = "Normal"
$ws.Range("E22").Value = '  +11.33%  '
$ws.Range("E23").Value = '  +0.19%  '
$c = $ws.Range("D24")
$c.Value = "'0.516"
$c.Style = "Normal"
$ws.Range("E24").Value = '  +14.18%  '
$c = $ws.Range("D25")
$c.Value = "'62.95"
$c.Style = "Normal"
$ws.Range("E25").Value = '  -1.79%  '
$c = $ws.Range("D26")
$c.Value = "'0.993"
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.54%  '
$ws.Range("E27").Value = '  -0.41%  '
$c = $ws.Range("D28")
$c.Value = "'7.70"
$c.Style = "Normal"
$ws.Range("E28").Value = '  +6.12%  '
$ws.Range("E29").Value = '  +1.22%  '
$ws.Range("E30").Value = '  -0.12%  '
$c = $ws.Range("D31")
$c.Value = "'6.27"
$c.Style = "Normal"
$ws.Range("E31").Value = '  +3.42%  '
$ws.Range("E32").Value = '  +1.23%  '
$c = $ws.Range("D33")
$c.Value = "'159.83"
$c.Style = "Normal"
$ws.Range("E33").Value = '  +1.58%  '
$c = $ws.Range("D34")
$c.Value = "'19.48"
$c.Style = "Normal"
$ws.Range("E34").Value = '  +2.36%  '
$c = $ws.Range("D35")
$c.Value = "'4.21"
$c.Style = "Normal"
$ws.Range("E35").Value = '  +4.56%  '
$c = $ws.Range("D36")
$c.Value = "'0.949"
$c.Style = "Normal"
$ws.Range("E36").Value = '  +7.18%  '
$c = $ws.Range("D37")
$c.Value = "'1.19"
$c.Style = "Normal"
$ws.Range("E37").Value = '  +4.84%  '
$ws.Range("E38").Value = '  +3.67%  '
$c = $ws.Range("D39")
$c.Value = "'37.74"
$c.Style = "Normal"
$ws.Range("E39").Value = '  +2.38%  '
$c = $ws.Range("D40")
$c.Value = "'0.856"
$c.Style = "Normal"
$ws.Range("E40").Value = '  -2.71%  '
$c = $ws.Range("D41")
$c.Value = "'3.77"
$c.Style = "Normal"
$ws.Range("E41").Value = '  +3.16%  '
$c = $ws.Range("D42")
$c.Value = "'299.44"
$c.Style = "Normal"
$ws.Range("E42").Value = '  +1.65%  '
$c = $ws.Range("D43")
$c.Value = "'141.59"
$c.Style = "Normal"
$ws.Range("E43").Value = '  +13.75%  '
$c = $ws.Range("D44")
$c.Value = "'0.996"
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range("D45")
$c.Value = "'0.607"
$c.Style = "Normal"
$ws.Range("E45").Value = '  +1.43%  '
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range("D46")
$c.Value = "'0.0979"
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.54%  '
$ws.Range("E47").Value = '  +4.14%  '
$c = $ws.Range("D48")
$c.Value = "'0.0542"
$c.Style = "Normal"
$ws.Range("E48").Value = '  +1.26%  '
$c = $ws.Range("D49")
$c.Value = "'10.68"
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range("D50")
$c.Value = "'19.39"
$c.Style = "Normal"
$ws.Range("E50").Value = '  +4.78%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D51")
$c.Value = "'4.79"
$c.Style = "Normal"
$ws.Range("E51").Value = '  +6.34%  '
